# Scheduled market-data refresh for Kujata_Profits (per-job-class leve profit sheets).
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H-N) with freshly scraped Universalis market data for the affected rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 13892887
$ws.Range("I62").Value = 18522600
$ws.Range("K62").Value = 18522600
$ws.Range("M62").Value = -18521976

# row 65
$ws.Range("H65").Value = 13892887
$ws.Range("I65").Value = 18522600
$ws.Range("K65").Value = 92613000
$ws.Range("M65").Value = -92609880

# row 129
$ws.Range("H129").Value = 659.5526
$ws.Range("I129").Value = 317.5
$ws.Range("J129").Value = 859.0833
$ws.Range("K129").Value = 952.5
$ws.Range("L129").Value = 2577.2499
$ws.Range("M129").Value = 4047.5
$ws.Range("N129").Value = -12577.2499

# row 135
$ws.Range("H135").Value = 22727752
$ws.Range("I135").Value = 219.1282
$ws.Range("K135").Value = 1972.1538
$ws.Range("M135").Value = 562.8462

# row 137
$ws.Range("H137").Value = 1341.5957
$ws.Range("I137").Value = 953.5
$ws.Range("K137").Value = 2860.5
$ws.Range("M137").Value = -310.5

# row 138
$ws.Range("H138").Value = 1247.9791
$ws.Range("I138").Value = 772.3415
$ws.Range("J138").Value = 1602.5454
$ws.Range("K138").Value = 2317.0245
$ws.Range("L138").Value = 4807.6362
$ws.Range("M138").Value = 2822.9755
$ws.Range("N138").Value = -15087.6362

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4553.623
$ws.Range("I32").Value = 4059.9836
$ws.Range("J32").Value = 8317.625
$ws.Range("K32").Value = 4059.9836
$ws.Range("L32").Value = 8317.625
$ws.Range("M32").Value = -3772.9836
$ws.Range("N32").Value = -8891.625

# row 61
$ws.Range("H61").Value = 100001944
$ws.Range("I61").Value = 125001670
$ws.Range("K61").Value = 125001670
$ws.Range("M61").Value = -125001458

# row 74
$ws.Range("H74").Value = 1736.6316
$ws.Range("I74").Value = 1291.6428
$ws.Range("K74").Value = 1291.6428
$ws.Range("M74").Value = -417.6428000000001

# row 77
$ws.Range("H77").Value = 1736.6316
$ws.Range("I77").Value = 1291.6428
$ws.Range("K77").Value = 6458.214
$ws.Range("M77").Value = -2090.214

# row 122
$ws.Range("H122").Value = 1421
$ws.Range("I122").Value = 1484.0625
$ws.Range("K122").Value = 4452.1875
$ws.Range("M122").Value = -2002.1875

# row 132
$ws.Range("H132").Value = 1463.7428
$ws.Range("I132").Value = 1407.0416
$ws.Range("J132").Value = 1587.4546
$ws.Range("K132").Value = 4221.1248
$ws.Range("L132").Value = 4762.3638
$ws.Range("M132").Value = -1691.1248
$ws.Range("N132").Value = -9822.363799999999

# row 136
$ws.Range("H136").Value = 100001944
$ws.Range("I136").Value = 125001670
$ws.Range("K136").Value = 375005010
$ws.Range("M136").Value = -375002460

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2014.931
$ws.Range("I31").Value = 1937.36
$ws.Range("J31").Value = 2499.75
$ws.Range("K31").Value = 1937.36
$ws.Range("L31").Value = 2499.75
$ws.Range("M31").Value = -1642.36
$ws.Range("N31").Value = -3089.75

# row 34
$ws.Range("H34").Value = 2014.931
$ws.Range("I34").Value = 1937.36
$ws.Range("J34").Value = 2499.75
$ws.Range("K34").Value = 1937.36
$ws.Range("L34").Value = 2499.75
$ws.Range("M34").Value = -1735.36
$ws.Range("N34").Value = -2903.75

# row 58
$ws.Range("H58").Value = 995.64703
$ws.Range("I58").Value = 923.48
$ws.Range("K58").Value = 923.48
$ws.Range("M58").Value = -720.48

# row 86
$ws.Range("H86").Value = 1862916.1
$ws.Range("I86").Value = 2905188.8
$ws.Range("J86").Value = 18895
$ws.Range("K86").Value = 2905188.8
$ws.Range("L86").Value = 18895
$ws.Range("M86").Value = -2904065.8
$ws.Range("N86").Value = -21141

# row 89
$ws.Range("H89").Value = 1862916.1
$ws.Range("I89").Value = 2905188.8
$ws.Range("J89").Value = 18895
$ws.Range("K89").Value = 14525944
$ws.Range("L89").Value = 94475
$ws.Range("M89").Value = -14520328
$ws.Range("N89").Value = -105707

# row 107
$ws.Range("H107").Value = 591.0952
$ws.Range("I107").Value = 488.84616
$ws.Range("J107").Value = 757.25
$ws.Range("K107").Value = 488.84616
$ws.Range("L107").Value = 757.25
$ws.Range("M107").Value = 1431.15384
$ws.Range("N107").Value = -4597.25

# row 132
$ws.Range("H132").Value = 3588.54
$ws.Range("I132").Value = 3988.0715
$ws.Range("K132").Value = 11964.2145
$ws.Range("M132").Value = -9434.2145

# row 134
$ws.Range("H134").Value = 11906109
$ws.Range("I134").Value = 1396.5405
$ws.Range("J134").Value = 100000980
$ws.Range("K134").Value = 4189.6215
$ws.Range("L134").Value = 300002940
$ws.Range("M134").Value = -1654.6215
$ws.Range("N134").Value = -300008010

# row 136
$ws.Range("H136").Value = 995.64703
$ws.Range("I136").Value = 923.48
$ws.Range("K136").Value = 2770.44
$ws.Range("M136").Value = -220.4400000000001

$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 22223538
$ws.Range("J131").Value = 1476.0541
$ws.Range("L131").Value = 4428.1623
$ws.Range("N131").Value = -14508.1623

$ws = $wb.Worksheets.Item("GSM")
# row 126
$ws.Range("H126").Value = 2870.7856
$ws.Range("J126").Value = 3673.875
$ws.Range("L126").Value = 11021.625
$ws.Range("N126").Value = -15961.625

# row 132
$ws.Range("H132").Value = 2365.0952
$ws.Range("I132").Value = 1979.1875
$ws.Range("K132").Value = 5937.5625
$ws.Range("M132").Value = -3407.5625

$ws = $wb.Worksheets.Item("LTW")
# row 100
$ws.Range("H100").Value = 1976.2222
$ws.Range("I100").Value = 1396.75
$ws.Range("K100").Value = 1396.75
$ws.Range("M100").Value = -855.75

# row 132
$ws.Range("H132").Value = 29626.25
$ws.Range("I132").Value = 1640.9565
$ws.Range("J132").Value = 79138.69500000001
$ws.Range("K132").Value = 4922.8695
$ws.Range("L132").Value = 237416.085
$ws.Range("M132").Value = -2392.8695
$ws.Range("N132").Value = -242476.085

# row 136
$ws.Range("H136").Value = 4577.5312
$ws.Range("I136").Value = 4802.893
$ws.Range("K136").Value = 14408.679
$ws.Range("M136").Value = -11858.679

$ws = $wb.Worksheets.Item("WVR")
# row 47
$ws.Range("N47").ClearContents()
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0

# row 58
$ws.Range("H58").Value = 15000
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15616

# row 132
$ws.Range("H132").Value = 3027
$ws.Range("I132").Value = 4590.0835
$ws.Range("J132").Value = 1321.8182
$ws.Range("K132").Value = 13770.2505
$ws.Range("L132").Value = 3965.4546
$ws.Range("M132").Value = -11240.2505
$ws.Range("N132").Value = -9025.454600000001

# row 136
$ws.Range("H136").Value = 584.3333
$ws.Range("J136").Value = 1550
$ws.Range("L136").Value = 4650
$ws.Range("N136").Value = -9750
